$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D:E (pushing batsman..sr from D:I to F:K)
$ws.Columns("D:E").Insert()

# Fill in the new header cells
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"

# Fill in the new data cells for row 2
$ws.Range("D2").Value = "Mumbai Indians"
$ws.Range("E2").Value = "Sunrisers Hyderabad"
